$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.76%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.62%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.690"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.00%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06108"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.00%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.657"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.49%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8501"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9214"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.54%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1398"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.93%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04760"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'17.27%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.17%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03076"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.78%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09053"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001533"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.50%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.71%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.05%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.150"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.30%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.96%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1304"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.095"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'6.09%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04235"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.31%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.63%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003806"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-19.23%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.89%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D40").Value = "'0.03858"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.53%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004079"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-34.82%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.01630"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'14.22%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002217"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.81%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005159"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.41%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'0.1355"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-43.74%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'36.18%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
